# Auto-generated edit script applying numeric corrections to the
# Leve profit calculation columns (H:N) across several sheets, per
# the scheduled pricing-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 1391.7273
$ws.Cells.Item(12, 9).Value = 163.75
$ws.Cells.Item(12, 11).Value = 163.75
$ws.Cells.Item(12, 13).Value = 6.25
$ws.Cells.Item(16, 8).Value = 12000
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 12000
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(16, 14).Value = -12460
$ws.Cells.Item(33, 8).Value = 371.27777
$ws.Cells.Item(33, 9).Value = 160.875
$ws.Cells.Item(33, 11).Value = 160.875
$ws.Cells.Item(33, 13).Value = 68.125
$ws.Cells.Item(70, 8).Value = 6387.657
$ws.Cells.Item(70, 10).Value = 6820.533
$ws.Cells.Item(70, 12).Value = 20461.599
$ws.Cells.Item(70, 14).Value = -21001.599
$ws.Cells.Item(73, 8).Value = 6387.657
$ws.Cells.Item(73, 10).Value = 6820.533
$ws.Cells.Item(73, 12).Value = 20461.599
$ws.Cells.Item(73, 14).Value = -22333.599
$ws.Cells.Item(98, 8).Value = 1302.6666
$ws.Cells.Item(98, 9).Value = 1302.6666
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 1302.6666
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 195.3334
$ws.Cells.Item(98, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 1302.6666
$ws.Cells.Item(122, 9).Value = 1302.6666
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 3907.9998
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -1457.9998
$ws.Cells.Item(122, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2035.64
$ws.Cells.Item(61, 9).Value = 1947.4762
$ws.Cells.Item(61, 11).Value = 1947.4762
$ws.Cells.Item(61, 13).Value = -1735.4762
$ws.Cells.Item(63, 8).Value = 2790.4285
$ws.Cells.Item(63, 9).Value = 1666.6666
$ws.Cells.Item(63, 10).Value = 3633.25
$ws.Cells.Item(63, 11).Value = 1666.6666
$ws.Cells.Item(63, 12).Value = 3633.25
$ws.Cells.Item(63, 13).Value = -980.6666
$ws.Cells.Item(63, 14).Value = -5005.25
$ws.Cells.Item(66, 8).Value = 2790.4285
$ws.Cells.Item(66, 9).Value = 1666.6666
$ws.Cells.Item(66, 10).Value = 3633.25
$ws.Cells.Item(66, 11).Value = 8333.333000000001
$ws.Cells.Item(66, 12).Value = 18166.25
$ws.Cells.Item(66, 13).Value = -4901.333000000001
$ws.Cells.Item(66, 14).Value = -25030.25
$ws.Cells.Item(75, 8).Value = 70659.664
$ws.Cells.Item(75, 9).Value = 70000
$ws.Cells.Item(75, 10).Value = 70989.5
$ws.Cells.Item(75, 11).Value = 70000
$ws.Cells.Item(75, 12).Value = 70989.5
$ws.Cells.Item(75, 13).Value = -69126
$ws.Cells.Item(75, 14).Value = -72737.5
$ws.Cells.Item(78, 8).Value = 70659.664
$ws.Cells.Item(78, 9).Value = 70000
$ws.Cells.Item(78, 10).Value = 70989.5
$ws.Cells.Item(78, 11).Value = 210000
$ws.Cells.Item(78, 12).Value = 212968.5
$ws.Cells.Item(78, 13).Value = -205632
$ws.Cells.Item(78, 14).Value = -221704.5
$ws.Cells.Item(102, 8).Value = 2181.5
$ws.Cells.Item(102, 9).Value = 2350.2856
$ws.Cells.Item(102, 11).Value = 2350.2856
$ws.Cells.Item(102, 13).Value = -728.2856000000002
$ws.Cells.Item(110, 8).Value = 7175.3335
$ws.Cells.Item(110, 9).Value = 6239.591
$ws.Cells.Item(110, 10).Value = 9748.625
$ws.Cells.Item(110, 11).Value = 6239.591
$ws.Cells.Item(110, 12).Value = 9748.625
$ws.Cells.Item(110, 13).Value = -4194.591
$ws.Cells.Item(110, 14).Value = -13838.625
$ws.Cells.Item(122, 8).Value = 3833.25
$ws.Cells.Item(122, 9).Value = 2944.4
$ws.Cells.Item(122, 11).Value = 8833.200000000001
$ws.Cells.Item(122, 13).Value = -6383.200000000001
$ws.Cells.Item(132, 8).Value = 2685.0833
$ws.Cells.Item(132, 9).Value = 1197.2632
$ws.Cells.Item(132, 11).Value = 3591.7896
$ws.Cells.Item(132, 13).Value = -1061.7896
$ws.Cells.Item(136, 8).Value = 2035.64
$ws.Cells.Item(136, 9).Value = 1947.4762
$ws.Cells.Item(136, 11).Value = 5842.4286
$ws.Cells.Item(136, 13).Value = -3292.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 8444.632
$ws.Cells.Item(105, 9).Value = 1921.25
$ws.Cells.Item(105, 11).Value = 1921.25
$ws.Cells.Item(105, 13).Value = -174.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 34926.883
$ws.Cells.Item(31, 9).Value = 1960
$ws.Cells.Item(31, 10).Value = 103857.63
$ws.Cells.Item(31, 11).Value = 1960
$ws.Cells.Item(31, 12).Value = 103857.63
$ws.Cells.Item(31, 13).Value = -1665
$ws.Cells.Item(31, 14).Value = -104447.63
$ws.Cells.Item(34, 8).Value = 34926.883
$ws.Cells.Item(34, 9).Value = 1960
$ws.Cells.Item(34, 10).Value = 103857.63
$ws.Cells.Item(34, 11).Value = 1960
$ws.Cells.Item(34, 12).Value = 103857.63
$ws.Cells.Item(34, 13).Value = -1758
$ws.Cells.Item(34, 14).Value = -104261.63
$ws.Cells.Item(58, 8).Value = 7199.095
$ws.Cells.Item(58, 9).Value = 4199
$ws.Cells.Item(58, 10).Value = 11199.223
$ws.Cells.Item(58, 11).Value = 4199
$ws.Cells.Item(58, 12).Value = 11199.223
$ws.Cells.Item(58, 13).Value = -3996
$ws.Cells.Item(58, 14).Value = -11605.223
$ws.Cells.Item(122, 8).Value = 7856.1333
$ws.Cells.Item(122, 9).Value = 3844.889
$ws.Cells.Item(122, 11).Value = 11534.667
$ws.Cells.Item(122, 13).Value = -9084.667000000001
$ws.Cells.Item(134, 8).Value = 2051.4412
$ws.Cells.Item(134, 9).Value = 1549.0646
$ws.Cells.Item(134, 11).Value = 4647.1938
$ws.Cells.Item(134, 13).Value = -2112.1938
$ws.Cells.Item(136, 8).Value = 7199.095
$ws.Cells.Item(136, 9).Value = 4199
$ws.Cells.Item(136, 10).Value = 11199.223
$ws.Cells.Item(136, 11).Value = 12597
$ws.Cells.Item(136, 12).Value = 33597.669
$ws.Cells.Item(136, 13).Value = -10047
$ws.Cells.Item(136, 14).Value = -38697.669

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 4740.6665
$ws.Cells.Item(22, 9).Value = 2110
$ws.Cells.Item(22, 11).Value = 6330
$ws.Cells.Item(22, 13).Value = -6161
$ws.Cells.Item(27, 8).Value = 4740.6665
$ws.Cells.Item(27, 9).Value = 2110
$ws.Cells.Item(27, 11).Value = 6330
$ws.Cells.Item(27, 13).Value = -6228
$ws.Cells.Item(39, 8).Value = 3083.3333
$ws.Cells.Item(39, 9).Value = 3750
$ws.Cells.Item(39, 10).Value = 2750
$ws.Cells.Item(39, 11).Value = 11250
$ws.Cells.Item(39, 12).Value = 8250
$ws.Cells.Item(39, 13).Value = -10956
$ws.Cells.Item(39, 14).Value = -8838
$ws.Cells.Item(133, 8).Value = 9217.25
$ws.Cells.Item(133, 9).Value = 8945.333000000001
$ws.Cells.Item(133, 11).Value = 26835.999
$ws.Cells.Item(133, 13).Value = -21775.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2293.4443
$ws.Cells.Item(102, 9).Value = 1218
$ws.Cells.Item(102, 11).Value = 1218
$ws.Cells.Item(102, 13).Value = 404
$ws.Cells.Item(122, 8).Value = 6023.385
$ws.Cells.Item(122, 9).Value = 4389.909
$ws.Cells.Item(122, 11).Value = 13169.727
$ws.Cells.Item(122, 13).Value = -10719.727
$ws.Cells.Item(126, 8).Value = 3849.25
$ws.Cells.Item(126, 9).Value = 1799
$ws.Cells.Item(126, 10).Value = 10000
$ws.Cells.Item(126, 11).Value = 5397
$ws.Cells.Item(126, 12).Value = 30000
$ws.Cells.Item(126, 13).Value = -2927
$ws.Cells.Item(126, 14).Value = -34940
$ws.Cells.Item(132, 8).Value = 63918.316
$ws.Cells.Item(132, 9).Value = 81817.07000000001
$ws.Cells.Item(132, 10).Value = 13801.8
$ws.Cells.Item(132, 11).Value = 245451.21
$ws.Cells.Item(132, 12).Value = 41405.39999999999
$ws.Cells.Item(132, 13).Value = -242921.21
$ws.Cells.Item(132, 14).Value = -46465.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 8209.223
$ws.Cells.Item(40, 9).Value = 8209.223
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 8209.223
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -8073.223
$ws.Cells.Item(40, 14).ClearContents()
$ws.Cells.Item(43, 8).Value = 16582.75
$ws.Cells.Item(43, 9).Value = 15545.272
$ws.Cells.Item(43, 11).Value = 15545.272
$ws.Cells.Item(43, 13).Value = -15352.272
$ws.Cells.Item(61, 8).Value = 12132
$ws.Cells.Item(61, 9).Value = 12127.286
$ws.Cells.Item(61, 10).Value = 12138.6
$ws.Cells.Item(61, 11).Value = 12127.286
$ws.Cells.Item(61, 12).Value = 12138.6
$ws.Cells.Item(61, 13).Value = -11925.286
$ws.Cells.Item(61, 14).Value = -12542.6
$ws.Cells.Item(113, 8).Value = 12132
$ws.Cells.Item(113, 9).Value = 12127.286
$ws.Cells.Item(113, 10).Value = 12138.6
$ws.Cells.Item(113, 11).Value = 12127.286
$ws.Cells.Item(113, 12).Value = 12138.6
$ws.Cells.Item(113, 13).Value = -9957.286
$ws.Cells.Item(113, 14).Value = -16478.6
$ws.Cells.Item(122, 8).Value = 4684.1816
$ws.Cells.Item(122, 9).Value = 3740.0688
$ws.Cells.Item(122, 10).Value = 11529
$ws.Cells.Item(122, 11).Value = 11220.2064
$ws.Cells.Item(122, 12).Value = 34587
$ws.Cells.Item(122, 13).Value = -8770.206399999999
$ws.Cells.Item(122, 14).Value = -39487
$ws.Cells.Item(132, 8).Value = 4618.1113
$ws.Cells.Item(132, 9).Value = 2648.2
$ws.Cells.Item(132, 11).Value = 7944.599999999999
$ws.Cells.Item(132, 13).Value = -5414.599999999999
$ws.Cells.Item(136, 8).Value = 4585.05
$ws.Cells.Item(136, 9).Value = 2046.0667
$ws.Cells.Item(136, 11).Value = 6138.2001
$ws.Cells.Item(136, 13).Value = -3588.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 561.5
$ws.Cells.Item(113, 9).Value = 549.6667
$ws.Cells.Item(113, 10).Value = 597
$ws.Cells.Item(113, 11).Value = 1649.0001
$ws.Cells.Item(113, 12).Value = 1791
$ws.Cells.Item(113, 13).Value = 520.9999
$ws.Cells.Item(113, 14).Value = -6131
$ws.Cells.Item(122, 8).Value = 12721.111
$ws.Cells.Item(122, 9).Value = 4394.4
$ws.Cells.Item(122, 10).Value = 15923.692
$ws.Cells.Item(122, 11).Value = 13183.2
$ws.Cells.Item(122, 12).Value = 47771.076
$ws.Cells.Item(122, 13).Value = -10733.2
$ws.Cells.Item(122, 14).Value = -52671.076
$ws.Cells.Item(126, 8).Value = 6570.278
$ws.Cells.Item(126, 9).Value = 6020.7856
$ws.Cells.Item(126, 11).Value = 18062.3568
$ws.Cells.Item(126, 13).Value = -15592.3568
$ws.Cells.Item(132, 8).Value = 4695.069
$ws.Cells.Item(132, 9).Value = 4955.1816
$ws.Cells.Item(132, 10).Value = 3877.5715
$ws.Cells.Item(132, 11).Value = 14865.5448
$ws.Cells.Item(132, 12).Value = 11632.7145
$ws.Cells.Item(132, 13).Value = -12335.5448
$ws.Cells.Item(132, 14).Value = -16692.7145
$ws.Cells.Item(135, 8).Value = 56690.453
$ws.Cells.Item(135, 10).Value = 56690.453
$ws.Cells.Item(135, 12).Value = 56690.453
$ws.Cells.Item(135, 14).Value = -66830.45300000001
